$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in newly entered grade values ---
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0

$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 5

$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 5

$ws.Range("J20").Value = 5

$ws.Range("H31").Value = 5
$ws.Range("I31").Value = 5

# --- New J column cells (need the same formatting as existing J cells, e.g. J20) ---
$ws.Range("J20").Copy() | Out-Null
$ws.Range("J4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("J4").Value = 0
$ws.Range("J12").Value = 5

# --- Totals column K: SUM formula for each student row ---
$ws.Range("K4").Formula = "=SUM(C4:J4)"
$ws.Range("K5:K32").Formula = "=SUM(C5:J5)"

# --- Conditional formatting: copy color scale from J4:J32 onto K4:K32 ---
$rangeK = $ws.Range("K4:K32")
$cs = $rangeK.FormatConditions.AddColorScale(3)

$rangeJ = $ws.Range("J4:J32")
$fcJ = $rangeJ.FormatConditions.Item(1)
$cs.Priority = 1
$fcJ.Priority = 2

# --- Update active selection to J12 ---
$ws.Range("J12").Select() | Out-Null

Write-Host "Edit applied successfully"
